# Aggiornamento dati fino al 21 marzo: aggiunge le righe 230-233
# (date seriali 44304-44307) in fondo alla tabella, replicando lo
# stile della colonna A (formato data) usato dalle righe precedenti.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @(44304, 2, 14, 226.4273006631085),
    @(44305, 1, 13, 210.2539220443151),
    @(44306, 3, 15, 242.600679281902),
    @(44307, 0, 14, 226.4273006631085)
)

$lastRow = 229
$r = $lastRow + 1
foreach ($row in $newRows) {
    # Copia la formattazione della colonna A dall'ultima riga esistente
    # cosi' la nuova cella data riceve lo stesso style (s="2").
    $ws.Range("A$lastRow").Copy()
    $ws.Range("A$r").PasteSpecial(-4122)

    $ws.Range("A$r").Value = $row[0]
    $ws.Range("B$r").Value = $row[1]
    $ws.Range("C$r").Value = $row[2]
    $ws.Range("D$r").Value = $row[3]

    $r = $r + 1
}
